$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header cells: "<Name>_old" -> "<Name>_FV2410" and "<Name>_new" -> "<Name>_FV2504"
$headerNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $headerNames.Length; $i++) {
    $colOld = 1 + $i        # columns A..J
    $colNew = 12 + $i       # columns L..U (K = "diff", unchanged)
    $ws.Cells.Item(1, $colOld).Value = "$($headerNames[$i])_FV2410"
    $ws.Cells.Item(1, $colNew).Value = "$($headerNames[$i])_FV2504"
}

# 2) Turn the data range into an Excel Table (ListObject), matching the target table1.xml.
#    Creating a table over cells that already carry bold/centered/bordered formatting makes
#    Excel record a header-row style "diff" (dxf); to avoid introducing that extra styling
#    info (not present in the target workbook) we temporarily park a copy of the header
#    formatting, strip the header formatting, build the table, then restore the formatting.
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A200:U200")

$headerRange.Copy($scratch)
$headerRange.ClearFormats()

$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U76"), [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"

$scratch.Copy()
$headerRange.PasteSpecial(-4122)  # xlPasteFormats
$scratch.Clear()
$excel.CutCopyMode = $false

# 3) Freeze the header row (top row) in the sheet view.
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
